# ---------------------------------------------------------------------------
# disability_prevalence.xlsx (Senaki) - update to "Unified database of
# targeted social assistance program" dataset:
#   - new title
#   - new "family with disabilities Persons" row of data
#   - "Number of disability persons" row renamed/retargeted to
#     "disabilities Persons" with new data
#   - source note row shifted down, now spans A:H on row 6
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at row 4 -- this pushes the old row4 (disability
#    numbers) to row5 and the old row5 (source note) to row6, keeping their
#    contents intact so we only have to touch what actually changed.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------------------
# 2. Text content
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Senaki Municipality"
$ws.Range("A2").Value = "(End of year, persons)"
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A5").Value = "disabilities Persons "
# A6 already holds the (unchanged) "Source: ..." rich text after the insert.

# ---------------------------------------------------------------------------
# 3. Numeric data
# ---------------------------------------------------------------------------
$row4vals = @(828,797,780,787,796,796,819,827)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $row4vals[$i]
}

$row5vals = @(938,908,881,885,885,888,917,922)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $row5vals[$i]
}

# Row6 (old row5) used to carry placeholder empty cells across B:H -- keep
# them empty (clear any leftover value from the old row).
foreach ($col in @("B","C","D","E","F","G","H")) {
    $ws.Range($col + "6").Value = $null
}

# ---------------------------------------------------------------------------
# 4. Merges -- title spans the whole table width now, source note shifted.
# ---------------------------------------------------------------------------
$ws.Range("A5:H5").UnMerge()
$ws.Range("A6:H6").Merge()
$ws.Range("A1:I1").Merge()

# ---------------------------------------------------------------------------
# 5. Row heights / column width / sheet formatting
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 51
$ws.Rows.Item(2).RowHeight = 14.5
$ws.Rows.Item(3).RowHeight = 18.75
$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 21
$ws.Rows.Item(6).RowHeight = 27.75

$ws.Columns.Item(1).ColumnWidth = 20.81640625

$ws.StandardHeight = 14.5

# ---------------------------------------------------------------------------
# 6. Styling
# ---------------------------------------------------------------------------

# -- column A default font (Georgian-capable "Sylfaen") --------------------
$colA = $ws.Columns.Item(1)
$colA.Font.Name = "Sylfaen"
$colA.Font.Size = 11
$colA.Font.Color = 0

# -- Row1 title: bold Arial 11, centered, wrapped, no border/fill ----------
$r1 = $ws.Range("A1:I1")
$r1.Font.Name = "Arial"
$r1.Font.Size = 11
$r1.Font.Bold = $true
$r1.Font.Underline = $false
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4108
$r1.WrapText = $true
$r1.Interior.Pattern = 0
foreach ($edge in 7,8,9,10) { $r1.Borders.Item($edge).LineStyle = -4142 }

# -- Row2: "(End of year, persons)" -----------------------------------------
$a2 = $ws.Range("A2")
$a2.Font.Name = "Arial"
$a2.Font.Size = 10
$a2.Font.Color = 0
$a2.Font.Bold = $false
$a2.HorizontalAlignment = -4131
$a2.VerticalAlignment = -4142
$a2.WrapText = $false
$a2.Interior.Pattern = 1
$a2.Interior.PatternColorIndex = -4105
$a2.Interior.ThemeColor = 4
$a2.Interior.TintAndShade = 0
foreach ($edge in 7,8,9,10) { $a2.Borders.Item($edge).LineStyle = -4142 }

# -- Row3: A3 blank corner cell (Sylfaen, top border only) ------------------
$a3 = $ws.Range("A3")
$a3.Font.Name = "Sylfaen"
$a3.Font.Size = 11
$a3.Font.Color = 0
$a3.Interior.Pattern = -4142
$a3.Borders.Item(8).LineStyle = 1
$a3.Borders.Item(8).Weight = 2
$a3.Borders.Item(9).LineStyle = -4142

# -- Row3: B3:I3 year headers -----------------------------------------------
$yearsRow = $ws.Range("B3:I3")
$yearsRow.Font.Name = "Arial"
$yearsRow.Font.Size = 10
$yearsRow.Font.Color = 0
$yearsRow.HorizontalAlignment = -4108
$yearsRow.VerticalAlignment = -4108
$yearsRow.WrapText = $true
$yearsRow.Interior.Pattern = 1
$yearsRow.Interior.PatternColorIndex = -4105
$yearsRow.Interior.ThemeColor = 4
$yearsRow.Interior.TintAndShade = 0
foreach ($edge in 8,9) {
    $yearsRow.Borders.Item($edge).LineStyle = 1
    $yearsRow.Borders.Item($edge).Weight = 2
}
foreach ($edge in 7,10) { $yearsRow.Borders.Item($edge).LineStyle = -4142 }

# -- Row4: A4 label ("family with disabilities Persons ") -------------------
$a4 = $ws.Range("A4")
$a4.Font.Name = "Arial"
$a4.Font.Size = 10
$a4.Font.Color = 0
$a4.HorizontalAlignment = -4131
$a4.VerticalAlignment = -4108
$a4.WrapText = $true
$a4.Interior.Pattern = 1
$a4.Interior.PatternColorIndex = -4105
$a4.Interior.ThemeColor = 4
$a4.Interior.TintAndShade = 0
$a4.Borders.Item(8).LineStyle = 1
$a4.Borders.Item(8).Weight = 2
$a4.Borders.Item(9).LineStyle = -4142

# -- Row4: B4:I4 numeric data -------------------------------------------------
$r4nums = $ws.Range("B4:I4")
$r4nums.Font.Name = "Arial"
$r4nums.Font.Size = 10
$r4nums.Font.Color = 0
$r4nums.NumberFormat = "#\ ##0"
$r4nums.HorizontalAlignment = -4131
$r4nums.Interior.Pattern = 1
$r4nums.Interior.PatternColorIndex = -4105
$r4nums.Interior.ThemeColor = 4
$r4nums.Interior.TintAndShade = 0
foreach ($edge in 7,8,9,10) { $r4nums.Borders.Item($edge).LineStyle = -4142 }

# -- Row5: A5 label ("disabilities Persons ") --------------------------------
$a5 = $ws.Range("A5")
$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.Font.Color = 0
$a5.HorizontalAlignment = -4131
$a5.VerticalAlignment = -4108
$a5.WrapText = $true
$a5.Interior.Pattern = 1
$a5.Interior.PatternColorIndex = -4105
$a5.Interior.ThemeColor = 4
$a5.Interior.TintAndShade = 0
$a5.Borders.Item(9).LineStyle = 1
$a5.Borders.Item(9).Weight = 2
$a5.Borders.Item(8).LineStyle = -4142

# -- Row5: B5:H5 numeric data -------------------------------------------------
$r5nums = $ws.Range("B5:H5")
$r5nums.Font.Name = "Arial"
$r5nums.Font.Size = 10
$r5nums.Font.Color = 0
$r5nums.NumberFormat = "#\ ##0"
$r5nums.HorizontalAlignment = -4131
$r5nums.Interior.Pattern = 1
$r5nums.Interior.PatternColorIndex = -4105
$r5nums.Interior.ThemeColor = 4
$r5nums.Interior.TintAndShade = 0
foreach ($edge in 7,8,9,10) { $r5nums.Borders.Item($edge).LineStyle = -4142 }

# -- Row5: I5 numeric data (bottom border, like the A5 label row end) -------
$i5 = $ws.Range("I5")
$i5.Font.Name = "Arial"
$i5.Font.Size = 10
$i5.Font.Color = 0
$i5.NumberFormat = "#\ ##0"
$i5.HorizontalAlignment = -4131
$i5.Interior.Pattern = 1
$i5.Interior.PatternColorIndex = -4105
$i5.Interior.ThemeColor = 4
$i5.Interior.TintAndShade = 0
$i5.Borders.Item(9).LineStyle = 1
$i5.Borders.Item(9).Weight = 2
foreach ($edge in 7,8,10) { $i5.Borders.Item($edge).LineStyle = -4142 }

# -- Row6: A6 source note + B6:H6 trailing merged cells ----------------------
$a6 = $ws.Range("A6")
$a6.Font.Name = "Arial"
$a6.Font.Size = 9
$a6.Font.Color = 0
$a6.Font.Bold = $false
$a6.Font.Underline = $false
$a6.HorizontalAlignment = -4131
$a6.VerticalAlignment = -4108
$a6.WrapText = $true
$a6.Interior.Pattern = 1
$a6.Interior.PatternColorIndex = -4105
$a6.Interior.ThemeColor = 4
$a6.Interior.TintAndShade = 0
foreach ($edge in 7,8,9,10) { $a6.Borders.Item($edge).LineStyle = -4142 }

$r6rest = $ws.Range("B6:H6")
$r6rest.Font.Name = "Arial"
$r6rest.Font.Size = 9
$r6rest.Font.Color = 0
$r6rest.HorizontalAlignment = -4131
$r6rest.VerticalAlignment = -4108
$r6rest.WrapText = $true
$r6rest.Interior.Pattern = 1
$r6rest.Interior.PatternColorIndex = -4105
$r6rest.Interior.ThemeColor = 4
$r6rest.Interior.TintAndShade = 0
$r6rest.Borders.Item(8).LineStyle = 1
$r6rest.Borders.Item(8).Weight = 2
foreach ($edge in 7,9,10) { $r6rest.Borders.Item($edge).LineStyle = -4142 }

# ---------------------------------------------------------------------------
# 7. Selection cosmetics (matches the saved workbook's last selection)
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Select()
